$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 currently holds the text "R40" (a shared string).
# The change replaces its content with the text "1" (kept as text,
# not a number, so the cell's underlying type stays the same).
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
